$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextCell 2 4 '23.660.41'
Set-TextCell 2 5 '  +1.82%  '

Set-TextCell 3 4 '1.657.01'

Set-TextCell 4 4 '0.9964'
Set-TextCell 4 5 '  -0.42%  '

Set-TextCell 5 4 '306.29'
Set-TextCell 5 5 '  +0.83%  '

Set-TextCell 6 4 '0.9971'
Set-TextCell 6 5 '  -0.34%  '

Set-TextCell 7 4 '0.3774'
Set-TextCell 7 5 '  +0.12%  '

Set-TextCell 8 4 '52.09'
Set-TextCell 8 5 '  +0.15%  '

Set-TextCell 9 4 '0.3666'
Set-TextCell 9 5 '  +0.92%  '

Set-TextCell 10 4 '1.270'
Set-TextCell 10 5 '  -0.53%  '

Set-TextCell 11 4 '0.08177'
Set-TextCell 11 5 '  +0.38%  '

Set-TextCell 12 4 '0.9965'
Set-TextCell 12 5 '  -0.39%  '

Set-TextCell 13 4 '23.22'
Set-TextCell 13 5 '  +1.71%  '

Set-TextCell 14 4 '6.706'
Set-TextCell 14 5 '  +1.90%  '

Set-TextCell 15 4 '0.00001280'
Set-TextCell 15 5 '  +2.27%  '

Set-TextCell 16 4 '7.400'
Set-TextCell 16 5 '  -0.20%  '

Set-TextCell 17 4 '1.653.02'
Set-TextCell 17 5 '  +3.22%  '

Set-TextCell 18 4 '95.49'
Set-TextCell 18 5 '  +1.52%  '

Set-TextCell 19 4 '0.06912'
Set-TextCell 19 5 '  -0.12%  '

Set-TextCell 20 4 '18.46'
Set-TextCell 20 5 '  +1.70%  '

Set-TextCell 21 5 '  +1.16%  '

Set-TextCell 22 4 '0.9967'
Set-TextCell 22 5 '  -0.68%  '

Set-TextCell 23 4 '23.655.87'
Set-TextCell 23 5 '  +1.77%  '

Set-TextCell 24 4 '12.95'
Set-TextCell 24 5 '  +0.41%  '

Set-TextCell 25 4 '3.153'
Set-TextCell 25 5 '  +3.60%  '

Set-TextCell 26 4 '2.410'
Set-TextCell 26 5 '  -1.72%  '

Set-TextCell 27 5 '  +0.97%  '

Set-TextCell 28 4 '151.05'
Set-TextCell 28 5 '  +0.75%  '

Set-TextCell 29 4 '5.330'
Set-TextCell 29 5 '  +1.11%  '

Set-TextCell 30 4 '137.02'
Set-TextCell 30 5 '  +1.03%  '

Set-TextCell 31 4 '2.327'
Set-TextCell 31 5 '  -2.27%  '

Set-TextCell 32 4 '1.837.31'
Set-TextCell 32 5 '  +3.28%  '

Set-TextCell 33 4 '6.911'
Set-TextCell 33 5 '  +2.57%  '

Set-TextCell 34 4 '11.08'
Set-TextCell 34 5 '  +7.18%  '

Set-TextCell 35 4 '0.9783'
Set-TextCell 35 5 '  +1.43%  '

Set-TextCell 36 4 '0.02877'
Set-TextCell 36 5 '  +4.86%  '

Set-TextCell 37 4 '6.386'
Set-TextCell 37 5 '  +4.24%  '

Set-TextCell 38 2 'Algorand'
Set-TextCell 38 3 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell 38 4 '0.2577'
Set-TextCell 38 5 '  +2.08%  '

Set-TextCell 39 2 'Hedera'
Set-TextCell 39 3 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell 39 4 '0.07411'
Set-TextCell 39 5 '  -1.01%  '

Set-TextCell 40 4 '0.08919'
Set-TextCell 40 5 '  +1.38%  '

Set-TextCell 41 4 '1.380'
Set-TextCell 41 5 '  -0.48%  '

Set-TextCell 42 4 '0.7176'
Set-TextCell 42 5 '  +1.08%  '

Set-TextCell 43 2 'EnergySwap'
Set-TextCell 43 3 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 43 4 '16.41'
Set-TextCell 43 5 '  +5.53%  '

Set-TextCell 44 2 'Aptos'
Set-TextCell 44 3 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell 44 4 '12.62'
Set-TextCell 44 5 '  +1.45%  '

Set-TextCell 45 4 '0.6651'
Set-TextCell 45 5 '  +1.65%  '

Set-TextCell 46 4 '2.373'
Set-TextCell 46 5 '  +2.45%  '

Set-TextCell 47 4 '4.029'
Set-TextCell 47 5 '  +0.55%  '

Set-TextCell 48 4 '0.9962'
Set-TextCell 48 5 '  -0.35%  '

Set-TextCell 49 4 '0.08041'
Set-TextCell 49 5 '  +1.31%  '

Set-TextCell 50 2 'Flow'
Set-TextCell 50 3 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
Set-TextCell 50 4 '1.231'
Set-TextCell 50 5 '  +1.97%  '

Set-TextCell 51 2 'Quant'
Set-TextCell 51 3 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell 51 4 '129.24'
Set-TextCell 51 5 '  -2.60%  '
